$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (D) and "Volume(1h)" (E) values for each data row (2-51),
# as refreshed by the GitHub Actions cron job.
$updates = @(
    @{ Row = 2;  D = '30.468.36';   E = '  -0.95%  ' }
    @{ Row = 3;  D = '2.090.31';    E = '  -1.23%  ' }
    @{ Row = 4;  D = '1.002';       E = '  +0.15%  ' }
    @{ Row = 5;  D = '329.21';      E = '  -0.82%  ' }
    @{ Row = 6;  D = $null;         E = '  +0.22%  ' }
    @{ Row = 7;  D = '0.5208';      E = '  +0.61%  ' }
    @{ Row = 8;  D = '0.4339';      E = '  -0.43%  ' }
    @{ Row = 9;  D = '51.67';       E = '  +13.36%  ' }
    @{ Row = 10; D = '0.08849';     E = '  -1.93%  ' }
    @{ Row = 11; D = '1.156';       E = '  -2.04%  ' }
    @{ Row = 12; D = '24.42';       E = '  -2.17%  ' }
    @{ Row = 13; D = '2.089.62';    E = '  -1.36%  ' }
    @{ Row = 14; D = '6.686';       E = '  -1.50%  ' }
    @{ Row = 15; D = '7.677';       E = '  +0.16%  ' }
    @{ Row = 16; D = '95.90';       E = '  -1.69%  ' }
    @{ Row = 17; D = '1.003';       E = '  +0.23%  ' }
    @{ Row = 18; D = '0.00001118';  E = '  -2.00%  ' }
    @{ Row = 19; D = '0.06581';     E = '  -0.55%  ' }
    @{ Row = 20; D = '19.19';       E = '  -0.10%  ' }
    @{ Row = 21; D = '1.001';       E = '  +0.25%  ' }
    @{ Row = 22; D = '6.277';       E = '  -2.50%  ' }
    @{ Row = 23; D = '30.505.06';   E = '  -1.50%  ' }
    @{ Row = 24; D = '12.16';       E = '  +1.45%  ' }
    @{ Row = 25; D = $null;         E = '  +2.81%  ' }
    @{ Row = 26; D = '2.330.76';    E = '  -1.49%  ' }
    @{ Row = 27; D = '22.22';       E = '  -3.44%  ' }
    @{ Row = 28; D = '2.585';       E = '  +0.92%  ' }
    @{ Row = 29; D = '162.04';      E = '  -1.23%  ' }
    @{ Row = 30; D = '131.37';      E = '  -2.03%  ' }
    @{ Row = 31; D = '1.190';       E = '  -0.23%  ' }
    @{ Row = 32; D = $null;         E = '  -0.30%  ' }
    @{ Row = 33; D = '1.660';       E = '  +7.38%  ' }
    @{ Row = 34; D = '6.138';       E = '  -1.76%  ' }
    @{ Row = 35; D = '3.895';       E = '  +0.13%  ' }
    @{ Row = 36; D = $null;         E = '  +5.59%  ' }
    @{ Row = 37; D = '0.02566';     E = '  -0.73%  ' }
    @{ Row = 38; D = '0.06802';     E = '  +0.21%  ' }
    @{ Row = 39; D = '5.456';       E = '  -2.76%  ' }
    @{ Row = 40; D = '12.64';       E = '  -0.20%  ' }
    @{ Row = 41; D = '0.2262';      E = '  +0.47%  ' }
    @{ Row = 42; D = '0.6909';      E = '  +1.62%  ' }
    @{ Row = 43; D = '1.263';       E = '  +0.41%  ' }
    @{ Row = 44; D = $null;         E = '  +0.35%  ' }
    @{ Row = 45; D = '0.6374';      E = '  +0.78%  ' }
    @{ Row = 46; D = '13.97';       E = '  -3.24%  ' }
    @{ Row = 47; D = '2.198';       E = '  -2.54%  ' }
    @{ Row = 48; D = '3.624';       E = '  -0.92%  ' }
    @{ Row = 49; D = '1.233';       E = '  +11.13%  ' }
    @{ Row = 50; D = $null;         E = '  -3.02%  ' }
    @{ Row = 51; D = '81.75';       E = '  -2.08%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Range("D$row")
        # Several "price" strings (e.g. "1.002", "0.5208") parse as plain
        # numbers. The source data must stay textual (as it was stored as
        # an inline string), so force Text formatting while writing it,
        # then restore the cell's original (default) style so no other
        # formatting is disturbed.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }

    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}
